# Single Pick sheet: add a "UPH" column (C) and flip SinglePickQuantity
# values (column B) from negative to positive.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SINGLE PICK")

# New header for column C
$ws.Range("C1").Value = "UPH"

# Row -> (positive SinglePickQuantity, UPH)
$data = @(
    @{ Row = 2;  B = 27;  C = 11.40845070422535 },
    @{ Row = 3;  B = 10;  C = 4.225352112676056 },
    @{ Row = 4;  B = 23;  C = 9.71830985915493 },
    @{ Row = 5;  B = 35;  C = 14.7887323943662 },
    @{ Row = 6;  B = 44;  C = 18.59154929577465 },
    @{ Row = 7;  B = 60;  C = 25.35211267605634 },
    @{ Row = 8;  B = 72;  C = 30.4225352112676 },
    @{ Row = 9;  B = 98;  C = 41.40845070422535 },
    @{ Row = 10; B = 50;  C = 21.12676056338028 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
